# Apply the weekly fruit/vegetable price update (Alcachofa, Mapocho Venta
# Directa de Santiago). The source data rows got reshuffled; this script
# writes the resulting values for the Fecha (D), Volumen (J), Precio
# minimo/maximo/promedio (K/L/M), Origen (O) and Precio $/Kg (P) columns
# for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44418
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 15000
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 500

# Row 4
$ws.Range("D4").Value = 44449
$ws.Range("J4").Value = 45

# Row 7
$ws.Range("D7").Value = 44467
$ws.Range("J7").Value = 35
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 12000
$ws.Range("P7").Value = 400

# Row 8
$ws.Range("D8").Value = 44453
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 12000
$ws.Range("P8").Value = 400

# Row 9
$ws.Range("D9").Value = 44446
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 14000
$ws.Range("P9").Value = 467

# Row 11
$ws.Range("D11").Value = 44474
$ws.Range("J11").Value = 45
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 10000
$ws.Range("P11").Value = 333

# Row 12
$ws.Range("D12").Value = 44421
$ws.Range("J12").Value = 25
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15400
$ws.Range("P12").Value = 513

# Row 13
$ws.Range("D13").Value = 44425
$ws.Range("J13").Value = 35
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 14000
$ws.Range("P13").Value = 467

# Row 14
$ws.Range("D14").Value = 44432
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 14000
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 467

$wb.Save()
